# Displaying sd available space
# Adds a new "Tiny" typography, switches several Translation rows to use it,
# and adds two new Translation rows for the SD-card available space label/value.

$wb = $excel.ActiveWorkbook

$wsTypography = $wb.Worksheets.Item("Typography")
$wsTranslation = $wb.Worksheets.Item("Translation")

# --- Typography sheet: add a new "Tiny" typography row (row 7) ---
$wsTypography.Range("B7").Value = "Tiny"
$wsTypography.Range("C7").Value = "segoeui.ttf"
$wsTypography.Range("D7").Value = 18
$wsTypography.Range("E7").Value = 2
$wsTypography.Range("F7").Value = "?"
$wsTypography.Range("G7").Value = ".kMGB"

# --- Translation sheet: switch existing rows from "Default" to "Tiny" typography ---
$wsTranslation.Range("C4").Value = "Tiny"
$wsTranslation.Range("C5").Value = "Tiny"
$wsTranslation.Range("C6").Value = "Tiny"
$wsTranslation.Range("D6").Value = "Center"
$wsTranslation.Range("C21").Value = "Tiny"
$wsTranslation.Range("C22").Value = "Tiny"
$wsTranslation.Range("C23").Value = "Tiny"
$wsTranslation.Range("D23").Value = "Center"

# --- Translation sheet: add new rows for SD available space ---
$wsTranslation.Range("B24").Value = "SingleUseId28"
$wsTranslation.Range("C24").Value = "Tiny"
$wsTranslation.Range("D24").Value = "Right"
$wsTranslation.Range("E24").Value = "LTR"
$wsTranslation.Range("F24").Value = "SD"

$wsTranslation.Range("B25").Value = "SingleUseId29"
$wsTranslation.Range("C25").Value = "Tiny"
$wsTranslation.Range("D25").Value = "Right"
$wsTranslation.Range("E25").Value = "LTR"
$wsTranslation.Range("F25").Value = "<value>"
